$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values re-pulled / recalculated for specific rows
$updates = @{
    7  = 1
    9  = 3
    13 = -3
    27 = -3
    29 = -3
    31 = -2
    33 = -3
    38 = 1
    42 = 1
    44 = -2
    51 = 1
    52 = -2
    54 = 0
    59 = -6
    63 = 3
    67 = -1
    70 = 1
    71 = 3
    72 = -1
    73 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
